$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "currency" column (C), which
# shifts currency/notionalPrincipal/.../description one column to the
# right (C->D, D->E, ... K->L). Formulas referencing the shifted columns
# (e.g. D7 =0.05*D2) are automatically rewritten by Excel to keep
# pointing at the same logical cells (E7 =0.05*E2).
$ws.Columns.Item(3).Insert()

# New column header
$ws.Range("C1").Value = "role"

# Mark the long (asset / income) positions vs. the short (expense)
# positions. The expense-type contracts (OES0001 - personnel expenses,
# OER0001 - rent expenses) are flagged "short" so they can be treated as
# negative cash flows; everything else is "long".
$ws.Range("C2:C9").Value = "long"
$ws.Range("C10:C11").Value = "short"

# Approximate the column width Excel would have auto-sized for the new
# column (matches the width used for the neighboring bestFit columns).
$ws.Columns.Item(3).ColumnWidth = 8.83

# Move the active cell selection to below the newly inserted column,
# mirroring where the author's cursor ended up after the edit.
$ws.Range("C12").Select()
